# Add a new worksheet "ODI Batting Extra" as the third (last) sheet and
# populate it with the extra ODI batting stats for match 4284.

$wb = $excel.ActiveWorkbook

# Grab the existing "ODI Batting" sheet (source of the header cell style
# we want the new sheet's header row to match) and the last sheet (so the
# new sheet gets appended after it, not inserted before the first one).
$sourceSheet = $wb.Worksheets.Item("ODI Batting")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# ---- Header row (bold / bordered / centered, matching the other sheets) ----
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")

$sourceSheet.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---- Data row ----
function Set-TextCell($cell, $value) {
    # Force text storage (matches t="inlineStr"/string cells in the target)
    # without leaving a lingering "Text" number-format style on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item(2, 1) "4284"
$ws.Cells.Item(2, 2).Value = 6
Set-TextCell $ws.Cells.Item(2, 3) "3"
Set-TextCell $ws.Cells.Item(2, 4) "1"
Set-TextCell $ws.Cells.Item(2, 5) "30.65%"
Set-TextCell $ws.Cells.Item(2, 6) "YES"

# Restore the original active sheet/selection (adding a sheet makes it
# active by default; the source workbook keeps "Player Info" selected).
$wb.Worksheets.Item(1).Activate()
